# Updated status/accomplishment data: remove the now-unused breakdown
# columns (No. of Sites/bldg Reverted, Not yet started, Under Procurement,
# On Going, Completed) and the DIFFERENCE column for the data rows,
# keeping the PREVIOUS ACCOMPLISHMENT (AL) column intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data rows 2-10: clear columns AB:AK (site/building breakdown counts)
$ws.Range("AB2:AK10").ClearContents()

# Data rows 2-10: clear column AM (DIFFERENCE)
$ws.Range("AM2:AM10").ClearContents()
